$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $t = $r.Text
    $trimmed = $t.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "7") {
        # Collapse to the end of the "7" text (before the paragraph mark)
        $end = $r.End - 1
        $insertRange = $d.Range($end, $end)
        $insertRange.InsertAfter("-8:30")
        break
    }
}
